$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookup")

# --- Header row: "Symptom Keywords" -> "Symptoms" ---
$ws.Range("I1").Value2 = "Symptoms"

# --- New vital-sign columns (B:G) for shock rows, plus refined symptom text ---
# Row 3 (Compensatory Shock) entered first, then row 2 (Decompensatory Shock),
# matching the authoring order captured in the workbook.
$ws.Range("B3").Value2 = "anxious, restless, disoriented"
$ws.Range("B2").Value2 = "worsening, eventually unresponsive"
$ws.Range("C2").Value2 = "Rapid, weak, eventually slows"
$ws.Range("E2").Value2 = "pale, cool, clammy"
$ws.Range("D2").Value2 = "Continues to increase, becomes shallower"
$ws.Range("F2").Value2 = "Falls, radial pulse weakens"
$ws.Range("G2").Value2 = "Slower to respond"
$ws.Range("C3").Value2 = "rapid, weak, eventually slows"
$ws.Range("D3").Value2 = "rapid and shallow"
$ws.Range("E3").Value2 = "pale, cool, and clammy"
$ws.Range("I3").Value2 = "nausea, vomiting, dizziness"

# --- Column widths: widen B/C, add a new D column width ---
$ws.Columns.Item(2).ColumnWidth = 29.333333333333336
$ws.Columns.Item(3).ColumnWidth = 23.666666666666668
$ws.Columns.Item(4).ColumnWidth = 16.333333333333336

# --- Sheet view: drop the stale scroll anchor, move selection to D17 ---
$null = $ws.Range("D17").Select()
